$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8

$ws.Range("G4").Value = 1.47
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 1.95
$ws.Range("K4").Value = 2.22
$ws.Range("L4").Value = 5.9
$ws.Range("N4").Value = 11.5
$ws.Range("U4").Value = 1.78
$ws.Range("V4").Value = 1.99
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 6.7
$ws.Range("Y4").Value = 6.7
$ws.Range("Z4").Value = 9.25
$ws.Range("AC4").Value = 11.25
$ws.Range("AD4").Value = 6.6
$ws.Range("AE4").Value = 13
$ws.Range("AH4").Value = 13.5
$ws.Range("AI4").Value = 32
$ws.Range("AJ4").Value = 15.5
$ws.Range("AK4").Value = 100
$ws.Range("AL4").Value = 50
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 6.9
$ws.Range("AP4").Value = 15
$ws.Range("AU4").Value = 7.3
$ws.Range("AW4").Value = 7.7
$ws.Range("AX4").Value = 37
$ws.Range("AZ4").Value = 250
$ws.Range("BA4").Value = 250
$ws.Range("BB4").Value = 450
